$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

# The old DatasetsTab query (currently in B2) moves down to become the
# ProjectsTab query in B3, and B2 receives the updated DatasetsTab query.
$oldQuery = $ws.Range("B2").Value()

$newQuery = "SELECT DISTINCT`n" +
    "    TRIM(REPLACE(ds.dataset_title, '  ', ' ')) AS ""Title"",`n" +
    "    ds.dataset_source_id AS ""Source ID"", `n" +
    "    ds.primary_disease AS ""Primary Disease"",`n" +
    "    -- CAST(ds.participant_count AS INT) AS ""Participants Count"",`n" +
    "    CAST(ds.sample_count AS INT) AS ""Sample Count""`n" +
    "FROM df_geo ds`n" +
    "ORDER BY CAST(ds.dataset_title AS TEXT) ASC;"

$ws.Range("B3").Value = $oldQuery
$ws.Range("B2").Value = $newQuery
